$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3, 4 and 5 (Id 76361837 / 82085799 / 82108936) were re-ordered so the
# three records are sorted by Startdatum instead of by Id: the record that
# used to sit in row 5 now sits in row 3, the old row-3 record moved down to
# row 4, and the old row-4 record moved down to row 5. Capture the original
# row contents first, then write them back in the rotated order.

$colFirst = 1   # A
$colLast  = 51  # AY

$row3 = $ws.Range($ws.Cells.Item(3, $colFirst), $ws.Cells.Item(3, $colLast)).Value2
$row4 = $ws.Range($ws.Cells.Item(4, $colFirst), $ws.Cells.Item(4, $colLast)).Value2
$row5 = $ws.Range($ws.Cells.Item(5, $colFirst), $ws.Cells.Item(5, $colLast)).Value2

# A few columns hold text that looks like a number/date ("Antal" column I,
# and the date columns Y/AA) - Excel's automatic type-detection would
# otherwise silently reinterpret them as numbers/dates when written back.
# Capture their original displayed text explicitly so the rotated values
# keep their original text representation.
$textCols = @(9, 25, 27)  # I, Y, AA

$row3Text = @{}
$row4Text = @{}
$row5Text = @{}
foreach ($c in $textCols) {
    $row3Text[$c] = $ws.Cells.Item(3, $c).Text
    $row4Text[$c] = $ws.Cells.Item(4, $c).Text
    $row5Text[$c] = $ws.Cells.Item(5, $c).Text
}

function Set-RotatedRow($destRow, $values, $textMap) {
    $destStart = $ws.Cells.Item($destRow, $colFirst)
    $destEnd = $ws.Cells.Item($destRow, $colLast)
    $ws.Range($destStart, $destEnd).Value2 = $values

    foreach ($c in $textCols) {
        $text = $textMap[$c]
        if ($text -ne "") {
            $cell = $ws.Cells.Item($destRow, $c)
            $cell.NumberFormat = "@"
            $cell.Value2 = $text
            $cell.Style = "Normal"
        }
    }
}

Set-RotatedRow 3 $row5 $row5Text
Set-RotatedRow 4 $row3 $row3Text
Set-RotatedRow 5 $row4 $row4Text
